# Generate Report for Handoff
#
# Flips the localization-status report from "Handed back" to
# "Ready for handoff" and refreshes the associated generate/handoff
# timestamps, on all three sheets (Overview, zh-cn, de-de). Also narrows
# the Status-ish columns to match the shorter text (mirrors Excel's
# column autosize after the text shrank).

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-29 09:03:47"

# --- zh-cn sheet ------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("H2").Value = "2016-08-29 09:03:43"

# --- de-de sheet --------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Ready for handoff"
$de.Range("H2").Value = "2016-08-29 09:03:47"

# --- Re-fit the now-narrower status columns --------------------------
# (closest attainable width to the recorded 17.216-char autofit result,
# given this engine's 1/6-character column-width grid)
$narrowWidth = 16.333333333333332
$ov.Columns.Item(5).ColumnWidth = $narrowWidth
$ov.Columns.Item(6).ColumnWidth = $narrowWidth
$zh.Columns.Item(3).ColumnWidth = $narrowWidth
$de.Columns.Item(3).ColumnWidth = $narrowWidth
